$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 30, shifting existing rows 30-95 down to 31-96
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with its data
$ws.Range("A30").Value = 5
$ws.Range("B30").Value = "Macroferia Regional de Talca"
$ws.Range("C30").Value = "Maule"
$ws.Range("D30").Value = 44487
$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 100112031
$ws.Range("G30").Value = "Poroto verde"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 100
$ws.Range("K30").Value = 40000
$ws.Range("L30").Value = 40000
$ws.Range("M30").Value = 40000
$ws.Range("N30").Value = "$/malla 25 kilos"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 1600
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"
